$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold font, border, centered alignment) from H1 onto the
# new header cells, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8
